$d = $word.ActiveDocument
$d.Content.Find.Execute("Primera entrega 25/6/2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Primera entrega 26/6/2019", 2)
